# Apply the updated cryptos list (Price/Volume(1h) columns) per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.613.43"
$ws.Range("E2").Value = "  -1.37%  "

$ws.Range("D3").Value = "2.515.82"
$ws.Range("E3").Value = "  -4.49%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'584.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.87%  "

$ws.Range("D6").Value = "'173.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.13%  "

$ws.Range("E8").Value = "  -1.76%  "

$ws.Range("D9").Value = "2.514.15"
$ws.Range("E9").Value = "  -4.54%  "

$ws.Range("E10").Value = "  -0.43%  "

$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("E12").Value = "  -3.44%  "

$ws.Range("D13").Value = "'5.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.27%  "

$ws.Range("E14").Value = "  -3.68%  "

$ws.Range("D15").Value = "2.961.95"
$ws.Range("E15").Value = "  -4.85%  "

$ws.Range("E16").Value = "  -2.87%  "

$ws.Range("D17").Value = "66.495.69"
$ws.Range("E17").Value = "  -1.52%  "

$ws.Range("D18").Value = "2.518.11"
$ws.Range("E18").Value = "  -4.10%  "

$ws.Range("D19").Value = "'7.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.56%  "

$ws.Range("E20").Value = "  -5.83%  "

$ws.Range("D21").Value = "'348.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.45%  "

$ws.Range("E22").Value = "  -2.47%  "

$ws.Range("E24").Value = "  +1.92%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").Value = "'69.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").Value = "'9.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.43%  "

$ws.Range("E28").Value = "  +0.23%  "

$ws.Range("E30").Value = "  -2.90%  "

$ws.Range("D31").Value = "'529.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.50%  "

$ws.Range("D32").Value = "'8.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.49%  "

$ws.Range("E33").Value = "  -2.48%  "

$ws.Range("E34").Value = "  -2.96%  "

$ws.Range("E35").Value = "  -3.95%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("D38").Value = "'155.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.99%  "

$ws.Range("D39").Value = "'18.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.17%  "

$ws.Range("E40").Value = "  +0.31%  "

$ws.Range("E41").Value = "  -2.85%  "

$ws.Range("E42").Value = "  -1.98%  "

$ws.Range("E43").Value = "  -2.19%  "

$ws.Range("E44").Value = "  +4.33%  "

$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("E46").Value = "  -1.21%  "

$ws.Range("D47").Value = "'148.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.18%  "

$ws.Range("E48").Value = "  -3.79%  "

$ws.Range("D50").Value = "'1.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.73%  "

$ws.Range("E51").Value = "  -8.45%  "
